# "Generate Report for Handback"
#
# The handback report records, for each target locale sheet, the file that
# was produced for hand-back (Latest Target File / Latest Handback File)
# together with the date/time the hand-back round finished, and flips the
# overall Status from "Ready for handoff" to "Handed back: in sync with
# en-US" everywhere that status is shown (Overview summary + each locale
# sheet's Status column).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdFile = "c1db27d1-0544-4b38-b01f-2158aa5e43f3.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04f847f5446fe4996e0a420fb9f8e6c33cb0164c/e2e/c1db27d1-0544-4b38-b01f-2158aa5e43f3.md"

# --- Overview sheet: update the per-locale status cells -------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText

# --- zh-cn sheet ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText

# Latest Target File (J2): hyperlink to the handed-back source doc
$wsZh.Hyperlinks.Add($wsZh.Range("J2"), $mdUrl, "", "", $mdFile)

# Latest Handback File (K2) / Latest Handback DateTime (L2)
$wsZh.Range("K2").Value = "c1db27d1-0544-4b38-b01f-2158aa5e43f3.33e0eaf2402ae7cbedf36e4fb5830ed86bcf4cc7.zh-cn.xlf"
$wsZh.Range("L2").Value = "2016-12-15 05:02:45"

# --- de-de sheet ------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText

# Latest Target File (J2): hyperlink to the handed-back source doc
$wsDe.Hyperlinks.Add($wsDe.Range("J2"), $mdUrl, "", "", $mdFile)

# Latest Handback File (K2) / Latest Handback DateTime (L2)
$wsDe.Range("K2").Value = "c1db27d1-0544-4b38-b01f-2158aa5e43f3.33e0eaf2402ae7cbedf36e4fb5830ed86bcf4cc7.de-de.xlf"
$wsDe.Range("L2").Value = "2016-12-15 05:03:05"

# --- Column widths: widen the Status / Latest Target File / Latest Handback
#     File columns so the longer text introduced above fits (matches the
#     autofit Excel performs once the cells hold longer strings). ----------
$wsOverview.Range("E1:F1").ColumnWidth = 29.9777050018311

$wsZh.Range("C1").ColumnWidth = 29.9777050018311
$wsZh.Range("J1:K1").ColumnWidth = 40

$wsDe.Range("C1").ColumnWidth = 29.9777050018311
$wsDe.Range("J1:K1").ColumnWidth = 40
